# Plasma_Gen_Transformer_Part List - add two new part-list rows (20-21)
# for the 1.25mm pitch right-angle 5-pin headers (SMT + Thru-Hole variants),
# matching the new structure introduced in the Transformer part list update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: 1.25mm Pitch Header, Surface Mount, Right Angle, 5-Pin
$ws.Range("C20").Value = "053261-0571"
$ws.Range("D20").Value = "Molex"
$ws.Range("D20").Borders.LineStyle = 1
$ws.Range("F20").Value = "1.25mm Pitch Header, Surface Mount, Right Angle, 5-Pin"

# Row 21: 1.25mm Pitch Header, Thru-Hole, Right Angle, 5-Pin
$ws.Range("C21").Value = "053048-0510"
$ws.Range("D21").Value = "Molex"
$ws.Range("D21").Borders.LineStyle = 1
$ws.Range("F21").Value = "1.25mm Pitch Header, Thru-Hole, Right Angle, 5-Pin"

# Leave the selection where the author left off while editing
[void]$ws.Range("P9").Select()
